$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of resale-number data (row 52), appended below the existing table.
# Date/Time/Weekday/Week columns are stored as text (matching the existing
# rows), the remaining columns are numeric.
$ws.Range("A52:D52").NumberFormat = "@"

$ws.Range("A52").Value = "2023-06-17"
$ws.Range("B52").Value = "11:35:14"
$ws.Range("C52").Value = "Saturday"
$ws.Range("D52").Value = "24"

$ws.Range("E52").Value = 121915
$ws.Range("F52").Value = 133265
$ws.Range("G52").Value = 161957
$ws.Range("H52").Value = 133293
$ws.Range("I52").Value = 177064
$ws.Range("J52").Value = 114583
$ws.Range("K52").Value = 201008
$ws.Range("L52").Value = 224680
$ws.Range("M52").Value = 175028
$ws.Range("N52").Value = 103285
$ws.Range("O52").Value = 39148
$ws.Range("P52").Value = 34014
$ws.Range("Q52").Value = 51728
$ws.Range("R52").Value = -1
$ws.Range("S52").Value = 36622
$ws.Range("T52").Value = -1
